# dodelani separatniho gui pro disky
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# disk_list: rows re-shuffled / edited (separate "disks" GUI data)
# ---------------------------------------------------------------------------
$disks = $wb.Worksheets.Item("disk_list")

# Row 1: 514_Teleflex -> 515_ZF
$disks.Range("A1").Value = "515_ZF"
$disks.Range("B1").Value = "Z"
$disks.Range("C1").Value = "\\10.9.250.100\08_Project_ZF_515\kamery"
$disks.Range("D1").Value = "jhvadmin"
$disks.Range("E1").Value = "jhvadm1n"

# Row 2: 515_ZF -> Domaci Nas
$disks.Range("A2").Value = "Domaci Nas"
$disks.Range("B2").Value = "S"
$disks.Range("C2").Value = "\\192.168.1.20\Data"
$disks.Range("D2").ClearContents()
$disks.Range("E2").ClearContents()

# Row 3: Domaci Nas -> 518_Valeo II
$disks.Range("A3").Value = "518_Valeo II"
$disks.Range("B3").Value = "VV2"
$disks.Range("C3").Value = "\\192.168.1.10\10_vision"
$disks.Range("D3").Value = "jhv_vision"
$disks.Range("E3").Value = "Jhv*2708"
$disks.Range("F3").Value = "Druha sít, ixon`nfj"

# Row 4: 518_Valeo II -> 474_B Austin
$disks.Range("A4").Value = "474_B Austin"
$disks.Range("B4").Value = "P"
$disks.Range("C4").Value = "\\10.96.205.166\DATA"
$disks.Range("E4").Value = "*Jhv2708"
$disks.Range("F4").Value = "10.96.205.166`t`nVisionNas_474B`t`n`t`t`t`t`t`tuser:JHV_Vision, omron `nPass:*Jhv2708h`nfhjgfds"

# Row 5: was 518_Valeo, now new entry 529_witt
$disks.Range("A5").Value = "529_witt"
$disks.Range("B5").Value = "Wj"
$disks.Range("C5").Value = "\\192.168.0.192\Dat"
$disks.Range("D5").Value = "Visio"
$disks.Range("E5").Value = "*Jhv270"
$disks.Range("F5").Value = "l"

# Row 6: was 474_B Austin, now 514_Teleflex (original row 1 content)
$disks.Range("A6").Value = "514_Teleflex"
$disks.Range("B6").Value = "T"
$disks.Range("C6").Value = "\\192.168.14.245\Data\Kamery"
$disks.Range("D6").Value = "Vision"
$disks.Range("F6").ClearContents()

# Row 7: was xfdx, now 518_Valeo (original row 5 content)
$disks.Range("A7").Value = "518_Valeo"
$disks.Range("B7").Value = "VV"
$disks.Range("C7").Value = "\\192.168.208.200\10_vision"
$disks.Range("D7").Value = "jhv_vision"
$disks.Range("E7").Value = "Jhv*2708"
$disks.Range("F7").Value = "první sít, ixon`n\\192.168.208.200\10_vision"

# Row 8: was 529_witte, now VUT
$disks.Range("A8").Value = "VUT"
$disks.Range("B8").Value = "V"
$disks.Range("C8").Value = "\\gigadisk2.ro.vutbr.cz\GIGADISK2\home\9\4\0\221049"
$disks.Range("D8").Value = "xhlava51@vutbr.cz"
$disks.Range("E8").ClearContents()

# ---------------------------------------------------------------------------
# Settings: default interface-selection setting changed 5 -> 0
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")
$settings.Range("B1").Value = 0

# ---------------------------------------------------------------------------
# projects_bin2 (hidden staging sheet): gains the two rows that were being
# edited on disk_list (their pre-VV2/pre-Wj snapshot values)
# ---------------------------------------------------------------------------
$bin2 = $wb.Worksheets.Item("projects_bin2")
$bin2.Range("A3").Value = "529_witt"
$bin2.Range("B3").Value = "Wj"
$bin2.Range("C3").Value = "\\192.168.0.192\Dat"
$bin2.Range("D3").Value = "Visio"
$bin2.Range("E3").Value = "*Jhv270"
$bin2.Range("F3").Value = "l"

$bin2.Range("A4").Value = "518_Valeo II"
$bin2.Range("B4").Value = "V"
$bin2.Range("C4").Value = "\\192.168.1.10\10_vision"
$bin2.Range("D4").Value = "jhv_vision"
$bin2.Range("E4").Value = "Jhv*2708"
$bin2.Range("F4").Value = "Druha sít, ixon`nfj"

# ---------------------------------------------------------------------------
# Settings_recources: default screenshot path + maximized-window flag
# ---------------------------------------------------------------------------
$res = $wb.Worksheets.Item("Settings_recources")
$res.Range("B3").Value = "C:/Users/kubah/Pictures/Screenshots/"
$res.Range("B22").Value = "ano"
